$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-22 Wednesday" "2025-01-23 Thursday"

Replace-Text "960÷2=480, 0" "152÷7=21, 5"
Replace-Text "743÷9=82, 5" "505÷9=56, 1"
Replace-Text "784÷7=112, 0" "992÷7=141, 5"
Replace-Text "890÷8=111, 2" "123÷6=20, 3"
Replace-Text "312÷7=44, 4" "219÷4=54, 3"

Replace-Text "657÷4=164, 1" "516÷8=64, 4"
Replace-Text "631÷4=157, 3" "234÷4=58, 2"
Replace-Text "880÷5=176, 0" "238÷9=26, 4"
Replace-Text "763÷8=95, 3" "667÷4=166, 3"
Replace-Text "755÷4=188, 3" "757÷7=108, 1"

Replace-Text "219÷8=27, 3" "283÷3=94, 1"
Replace-Text "396÷4=99, 0" "103÷9=11, 4"
Replace-Text "245÷5=49, 0" "733÷2=366, 1"
Replace-Text "569÷8=71, 1" "778÷3=259, 1"
Replace-Text "882÷5=176, 2" "865÷5=173, 0"

Replace-Text "574÷3=191, 1" "943÷9=104, 7"
Replace-Text "550÷4=137, 2" "587÷3=195, 2"
Replace-Text "223÷4=55, 3" "290÷2=145, 0"
Replace-Text "553÷7=79, 0" "461÷9=51, 2"
Replace-Text "980÷9=108, 8" "692÷3=230, 2"

Replace-Text "331÷8=41, 3" "106÷5=21, 1"
Replace-Text "293÷3=97, 2" "397÷6=66, 1"
Replace-Text "685÷9=76, 1" "839÷3=279, 2"
Replace-Text "139÷5=27, 4" "643÷7=91, 6"
Replace-Text "590÷8=73, 6" "778÷5=155, 3"
